$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("航天发展", "航天发展", "航天发展"),
    @("三花智控", "三花智控", "平潭发展"),
    @("雷科防务", "雷科防务", "三花智控"),
    @("中国卫星", "中国卫星", "万向钱潮"),
    @("神剑股份", "平潭发展", "天普股份"),
    @("翠微股份", "山子高科", "山子高科"),
    @("五洲新春", "神剑股份", "神剑股份"),
    @("平潭发展", "翠微股份", "中超控股"),
    @("山子高科", "金风科技", "雷科防务"),
    @("万向钱潮", "五洲新春", "泰尔股份"),
    @("中超控股", "中超控股", "福龙马"),
    @("天奇股份", "万向钱潮", "天奇股份"),
    @("中国卫通", "中国卫通", "中国卫星"),
    @("锋龙股份", "卧龙电驱", "五洲新春"),
    @("金风科技", "昊志机电", "东百集团"),
    @("福龙马", "福龙马", "翠微股份"),
    @("昊志机电", "华菱线缆", "首开股份"),
    @("泰尔股份", "泰尔股份", "卧龙电驱"),
    @("卧龙电驱", "天奇股份", "大业股份"),
    @("蓝色光标", "拓维信息", "华菱线缆")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][2]
}
